$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet view: zoom to 55% and change selection to R27
$ws.Application.ActiveWindow.Zoom = 55
$ws.Range("R27").Select()

# Coordinate cleanups (Latitude = column D, Longitude = column E)
$ws.Range("E2").Value = 4.7534000000000001

$ws.Range("D3").Value = 52.3508
$ws.Range("E3").Value = 5.2647000000000004

$ws.Range("D4").Value = 52.111199999999997
$ws.Range("E4").Value = 4.6473000000000004

$ws.Range("D7").Value = 52.011600000000001
$ws.Range("E7").Value = 4.3571

$ws.Range("E9").Value = 4.6901000000000002

$ws.Range("D10").Value = 52.011499999999998
$ws.Range("E10").Value = 4.7104999999999997

$ws.Range("E11").Value = 4.646299

$ws.Range("D13").Value = 52.163600000000002
$ws.Range("E13").Value = 4.4802
$ws.Range("X13").Value = 0

$ws.Range("D14").Value = 51.924399999999999
$ws.Range("E14").Value = 4.4776999999999996

$ws.Range("D15").Value = 52.300400000000003
$ws.Range("E15").Value = 4.6744000000000003

$ws.Range("D16").Value = 51.856200000000001
$ws.Range("E16").Value = 4.2972000000000001

$ws.Range("E18").Value = 4.3494000000000002

$ws.Range("D20").Value = 52.079799999999999
$ws.Range("E20").Value = 4.8627000000000002
$ws.Range("Q20").Value = 0

$ws.Range("D21").Value = 52.442
$ws.Range("E21").Value = 4.8292000000000002

$ws.Range("D22").Value = 52.060699999999997
$ws.Range("E22").Value = 4.4939999999999998
